# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) figures to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 4815
$ws1.Range("F7").Value  = 118
$ws1.Range("F8").Value  = 0
$ws1.Range("F9").Value  = 95
$ws1.Range("F10").Value = 760
$ws1.Range("F12").Value = 1185
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 85
$ws1.Range("F19").Value = 0
$ws1.Range("F21").Value = 0
$ws1.Range("F22").Value = 39
$ws1.Range("F23").Value = 87
$ws1.Range("F28").Value = 46
$ws1.Range("F33").Value = 143
$ws1.Range("F34").Value = 289
$ws1.Range("F37").Value = 0
$ws1.Range("F38").Value = 9
$ws1.Range("F40").Value = 0
$ws1.Range("F42").Value = 71
$ws1.Range("F43").Value = 60
$ws1.Range("F48").Value = 587

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 33
$ws4.Range("F3").Value  = 230
$ws4.Range("F4").Value  = 4815
$ws4.Range("F7").Value  = 118
$ws4.Range("F8").Value  = 0
$ws4.Range("F9").Value  = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 225
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 112
$ws4.Range("F15").Value = 259
$ws4.Range("F17").Value = 85
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 112
$ws4.Range("F20").Value = 4009
$ws4.Range("F21").Value = 6342
$ws4.Range("F22").Value = 39
$ws4.Range("F23").Value = 39
$ws4.Range("F25").Value = 541
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 407
$ws4.Range("F30").Value = 20
$ws4.Range("F33").Value = 0
$ws4.Range("F34").Value = 143
$ws4.Range("F36").Value = 309
$ws4.Range("F37").Value = 374
$ws4.Range("F38").Value = 177
$ws4.Range("F40").Value = 1564
$ws4.Range("F41").Value = 0
$ws4.Range("F45").Value = 0
$ws4.Range("F47").Value = 2
$ws4.Range("F49").Value = 0
